$wb = $excel.ActiveWorkbook

# Calibrate capacity addition and retirement parameters on the
# "CRtPaL-profits" sheet: set the response coefficients (B2:B15) to 2.
$ws = $wb.Worksheets.Item("CRtPaL-profits")
$ws.Range("B2:B15").Value = 2

# Mirror the author's on-screen selection at save time.
$ws.Range("B2:B15").Select()

# Restore the originally active tab ("About") so the workbook-level
# active-sheet/tab state is unchanged, matching the source edit.
$wb.Worksheets.Item("About").Activate()
